$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 6 (pushes the CLOSEBROWSER row down to row 8)
$ws.Rows.Item(6).Resize(2).Insert()

# Copy formatting from row 5 (GOTOURL-style data row) into the two new rows (6 and 7)
$ws.Range("A5:E5").Copy()
$ws.Range("A6:E7").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

# Row 5 changes: SORTBY -> SCROLLINTOVIEW, sortOption -> sortByLabel, clear value
$ws.Range("B5").Value = "SCROLLINTOVIEW"
$ws.Range("C5").Value = "sortByLabel"
$ws.Range("E5").Value = ""

# Row 6 (new): CHOOSEFROM / sortOption / css / Price: Lowest first
$ws.Range("A6").Value = "TC # 01.01"
$ws.Range("B6").Value = "CHOOSEFROM"
$ws.Range("C6").Value = "sortOption"
$ws.Range("D6").Value = "css"
$ws.Range("E6").Value = "Price: Lowest first"

# Row 7 (new): VERIFYLOWESTPRICEFIRST / pageItems / css / blank
$ws.Range("A7").Value = "TC # 01.01"
$ws.Range("B7").Value = "VERIFYLOWESTPRICEFIRST"
$ws.Range("C7").Value = "pageItems"
$ws.Range("D7").Value = "css"
$ws.Range("E7").Value = ""

# Widen column B slightly (closest achievable width to the authored 23.21875)
$ws.Columns.Item(2).ColumnWidth = 22.3

# Update active selection
$ws.Range("C17").Select() | Out-Null
